# A new batch of NCM log entries was captured. The sheet is an append/shift
# log: the 5 freshest rows are written into rows 2-6 (in place, overwriting
# what used to be the 2 oldest rows there), the surviving older rows shuffle
# down by the same in-place-overwrite pattern, and the 3 rows that no longer
# fit in the old range are appended as brand-new rows 10-12. Column D is
# also widened to fit the longer description text now stored there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NcmRow([int]$row, [string]$id, [string]$ipi, [string]$ncm, [string]$desc, [string]$data) {
    $ws.Cells.Item($row, 1).Value = $id
    if ($ipi -ne "") {
        $ws.Cells.Item($row, 2).Value = $ipi
    } else {
        $ws.Cells.Item($row, 2).ClearContents()
    }
    $c = $ws.Cells.Item($row, 3)
    # Force NCM codes that look numeric ("01.01", "1", ...) to stay text,
    # matching every other cell in the column.
    $c.NumberFormat = "@"
    $c.Value = $ncm
    $ws.Cells.Item($row, 4).Value = $desc
    $ws.Cells.Item($row, 5).Value = $data
}

Set-NcmRow 2  "30005f7a-2b93-4ab4-ba24-cebf5e09ff8c" "NT" "0106.14.00" "Coelhos e lebres" "2025-05-21 14:06:26.557000"
Set-NcmRow 3  "83c927d9-9d56-4428-82eb-f69cfdf448db" ""   "01.01"      "CAVALOS, ASININOS E MUARES, VIVOS." "2025-05-21 13:59:23.436000"
Set-NcmRow 4  "9b5af8eb-b5df-47bd-9cba-3d798003a9fd" ""   "1"          "ANIMAIS VIVOS" "2025-05-21 13:59:22.588000"
Set-NcmRow 5  "911823be-8648-4824-89c1-ceb6e50a8666" "NT" "0106.12.00" "Baleias, golfinhos e botos (mamíferos da ordem Cetacea); peixes-boi (manatins) e dugongos (mamíferos da ordem Sirenia); otárias e focas, leões-marinhos e morsas (mamíferos da subordem Pinnipedia)" "2025-05-21 13:48:36.378000"
Set-NcmRow 6  "4690ebbf-6196-4d85-848c-d50b750fe6cd" "NT" "0105.12.00" "Peruas e perus" "2025-05-21 13:48:35.319000"
Set-NcmRow 7  "e1264496-3d8a-4650-8ee2-829419786d9e" ""   "01.05"      "AVES DA ESPÉCIE GALLUS DOMESTICUS, PATOS, GANSOS, PERUS, PERUAS E GALINHAS-D’ANGOLA (PINTADAS), DAS ESPÉCIES DOMÉSTICAS, VIVOS." "2025-05-20 19:06:22.392000"
Set-NcmRow 8  "3831aba8-4ff2-4c7c-859c-96789b1675c7" ""   "01.05"      "AVES DA ESPÉCIE GALLUS DOMESTICUS, PATOS, GANSOS, PERUS, PERUAS E GALINHAS-D’ANGOLA (PINTADAS), DAS ESPÉCIES DOMÉSTICAS, VIVOS." "2025-05-20 18:57:28.911000"
Set-NcmRow 9  "f768e507-e1c8-4764-ad06-cfc6ff08999a" ""   "0102.39.1"  "PARA REPRODUÇÃO" "2025-05-20 18:57:28.495000"
Set-NcmRow 10 "105f25df-4ffc-4487-b3bb-aae17a07d567" ""   "0102.39.1"  "PARA REPRODUÇÃO" "2025-05-20 18:53:10.311000"
Set-NcmRow 11 "e58092ef-e69a-47ec-9b6c-67245eac6716" ""   "1"          "ANIMAIS VIVOS" "2025-05-20 18:41:05.422000"
Set-NcmRow 12 "c1d3f1f1-0422-4883-a063-75d18b4d26d1" ""   "1"          "ANIMAIS VIVOS" "2025-05-20 16:16:30.735000"

# Widen the description column to fit the new, longer text. Excel stores
# column widths in "characters" (Maximum Digit Width units) with a small
# fixed pixel padding baked in, so ColumnWidth=197 would actually be
# serialised as 196 characters; 196.1667 is what round-trips to the target
# stored width of 197.
$ws.Columns("D").ColumnWidth = 196.1667
